$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as literal TEXT (avoids Excel's
# automatic "looks like a number" conversion) without leaving any
# numberFormat / quotePrefix residue on the cell's style.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '=""&"' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# -----------------------------------------------------------------
# 1) Insert two new columns at F:G. This pushes the "DataSunt Alb_1"
#    box (H:I) to J:K, and the "DataSunt Alb_2" box (L:M) to N:O,
#    while leaving the B:E box, its merges and its styles untouched.
# -----------------------------------------------------------------
$ws.Range("F1:G1").EntireColumn.Insert()

# -----------------------------------------------------------------
# 2) Turn the old "E" (right edge) column of the first box into a
#    "middle" column, and make the new "F" column the right edge,
#    adding the "Trama_Completa" header. Clear the spare "G" column
#    that the insert created.
# -----------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G2:G3").Clear()

$ws.Range("B2:F2").MergeCells = $true

$ws.Range("B3:F3").HorizontalAlignment = -4108  # xlCenter

Set-TextValue $ws.Range("B3") "Codigo_Caja"
Set-TextValue $ws.Range("C3") "Hojas_Totales"
Set-TextValue $ws.Range("D3") "Hojas_Leer"
Set-TextValue $ws.Range("E3") "NºMáquina"
Set-TextValue $ws.Range("F3") "Trama_Completa"

# -----------------------------------------------------------------
# 3) Add the data rows 4-12 (rows 4-5 underlined, rows 6-12 plain)
# -----------------------------------------------------------------
$codigoCaja = "184123450000217702"
$hojasTot   = "2"
$hojasLeer  = "1"
$trama      = "184123450000217702020201d0"

for ($r = 4; $r -le 5; $r++) {
    Set-TextValue $ws.Range("B$r") $codigoCaja
    Set-TextValue $ws.Range("C$r") $hojasTot
    Set-TextValue $ws.Range("D$r") $hojasTot
    Set-TextValue $ws.Range("E$r") $hojasLeer
    Set-TextValue $ws.Range("F$r") $trama
    $ws.Range("B$r`:E$r").Font.Underline = $true
}

for ($r = 6; $r -le 12; $r++) {
    Set-TextValue $ws.Range("B$r") $codigoCaja
    Set-TextValue $ws.Range("C$r") $hojasTot
    Set-TextValue $ws.Range("D$r") $hojasTot
    Set-TextValue $ws.Range("E$r") $hojasLeer
    Set-TextValue $ws.Range("F$r") $trama
}

# -----------------------------------------------------------------
# 4) Misc view tweaks
# -----------------------------------------------------------------
$ws.Range("G13").Select()
